$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 91.34108066666666
$ws.Range("H2").Value = 274.023242
$ws.Range("I2").Value = 0.2190334467302001
$ws.Range("J2").Value = 0.2190334467302
$ws.Range("M2").Value = 1.660421
$ws.Range("N2").Value = 4.981262999999999
$ws.Range("O2").Value = 0.03714789785507311
$ws.Range("P2").Value = 0.03714789785507311
$ws.Range("Q2").Value = 151.6646485016273
$ws.Range("R2").Value = 1364.981836514646
$ws.Range("S2").Value = 0.00813663210597807
$ws.Range("T2").Value = 0.008136632105978068

# Row 3
$ws.Range("G3").Value = 91.34108066666666
$ws.Range("H3").Value = 274.023242
$ws.Range("I3").Value = 0.2190334467302001
$ws.Range("J3").Value = 0.2190334467302
$ws.Range("O3").Value = 0.5631392661118858
$ws.Range("P3").Value = 0.5631392661118859
$ws.Range("Q3").Value = 2299.142718264466
$ws.Range("R3").Value = 20692.2844643802
$ws.Range("S3").Value = 0.1233463344456017
$ws.Range("T3").Value = 0.1233463344456017

# Row 4
$ws.Range("G4").Value = 91.34108066666666
$ws.Range("H4").Value = 274.023242
$ws.Range("I4").Value = 0.2190334467302001
$ws.Range("J4").Value = 0.2190334467302
$ws.Range("M4").Value = 17.866195
$ws.Range("N4").Value = 53.598585
$ws.Range("O4").Value = 0.399712836033041
$ws.Range("P4").Value = 0.399712836033041
$ws.Range("Q4").Value = 1631.917558701397
$ws.Range("R4").Value = 14687.25802831257
$ws.Range("S4").Value = 0.08755048017862029
$ws.Range("T4").Value = 0.08755048017862027

# Row 5
$ws.Range("G5").Value = 276.4348856666666
$ws.Range("H5").Value = 829.3046569999999
$ws.Range("I5").Value = 0.6628833966285105
$ws.Range("J5").Value = 0.6628833966285105
$ws.Range("M5").Value = 1.660421
$ws.Range("N5").Value = 4.981262999999999
$ws.Range("O5").Value = 0.03714789785507311
$ws.Range("P5").Value = 0.03714789785507311
$ws.Range("Q5").Value = 458.9982892935322
$ws.Range("R5").Value = 4130.98460364179
$ws.Range("S5").Value = 0.02462472470777982
$ws.Range("T5").Value = 0.02462472470777982

# Row 6
$ws.Range("G6").Value = 276.4348856666666
$ws.Range("H6").Value = 829.3046569999999
$ws.Range("I6").Value = 0.6628833966285105
$ws.Range("J6").Value = 0.6628833966285105
$ws.Range("O6").Value = 0.5631392661118858
$ws.Range("P6").Value = 0.5631392661118859
$ws.Range("Q6").Value = 6958.131541865199
$ws.Range("R6").Value = 62623.1838767868
$ws.Range("S6").Value = 0.3732956694951335
$ws.Range("T6").Value = 0.3732956694951336

# Row 7
$ws.Range("G7").Value = 276.4348856666666
$ws.Range("H7").Value = 829.3046569999999
$ws.Range("I7").Value = 0.6628833966285105
$ws.Range("J7").Value = 0.6628833966285105
$ws.Range("M7").Value = 17.866195
$ws.Range("N7").Value = 53.598585
$ws.Range("O7").Value = 0.399712836033041
$ws.Range("P7").Value = 0.399712836033041
$ws.Range("Q7").Value = 4938.839572123371
$ws.Range("R7").Value = 44449.55614911034
$ws.Range("S7").Value = 0.2649630024255971
$ws.Range("T7").Value = 0.2649630024255971

# Row 8
$ws.Range("G8").Value = 49.24290466666667
$ws.Range("H8").Value = 147.728714
$ws.Range("I8").Value = 0.1180831566412894
$ws.Range("J8").Value = 0.1180831566412894
$ws.Range("M8").Value = 1.660421
$ws.Range("N8").Value = 4.981262999999999
$ws.Range("O8").Value = 0.03714789785507311
$ws.Range("P8").Value = 0.03714789785507311
$ws.Range("Q8").Value = 81.76395300953132
$ws.Range("R8").Value = 735.8755770857819
$ws.Range("S8").Value = 0.004386541041315218
$ws.Range("T8").Value = 0.004386541041315217

# Row 9
$ws.Range("G9").Value = 49.24290466666667
$ws.Range("H9").Value = 147.728714
$ws.Range("I9").Value = 0.1180831566412894
$ws.Range("J9").Value = 0.1180831566412894
$ws.Range("O9").Value = 0.5631392661118858
$ws.Range("P9").Value = 0.5631392661118859
$ws.Range("Q9").Value = 1239.491200062781
$ws.Range("R9").Value = 11155.42080056503
$ws.Range("S9").Value = 0.0664972621711506
$ws.Range("T9").Value = 0.0664972621711506

# Row 10
$ws.Range("G10").Value = 49.24290466666667
$ws.Range("H10").Value = 147.728714
$ws.Range("I10").Value = 0.1180831566412894
$ws.Range("J10").Value = 0.1180831566412894
$ws.Range("M10").Value = 17.866195
$ws.Range("N10").Value = 53.598585
$ws.Range("O10").Value = 0.399712836033041
$ws.Range("P10").Value = 0.399712836033041
$ws.Range("Q10").Value = 879.7833371410768
$ws.Range("R10").Value = 7918.05003426969
$ws.Range("S10").Value = 0.04719935342882363
$ws.Range("T10").Value = 0.04719935342882361

